$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.362.91'
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = '1.846.47'

$ws.Range('D4').Value = '''0.9986'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '''240.53'
$ws.Range('E5').Value = '  -0.51%  '

$ws.Range('D6').Value = '''0.6305'
$ws.Range('E6').Value = '  +0.34%  '

$ws.Range('D7').Value = '''0.9999'
$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '''0.07537'
$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('D9').Value = '''0.2956'
$ws.Range('E9').Value = '  -0.70%  '

$ws.Range('D10').Value = '''24.42'
$ws.Range('E10').Value = '  +0.22%  '

$ws.Range('D11').Value = '''0.07700'
$ws.Range('E11').Value = '  -0.24%  '

$ws.Range('D12').Value = '1.857.09'
$ws.Range('E12').Value = '  -2.74%  '

$ws.Range('D13').Value = '''4.993'

$ws.Range('D14').Value = '''0.6830'
$ws.Range('E14').Value = '  -0.94%  '

$ws.Range('E15').Value = '  +2.19%  '

$ws.Range('D16').Value = '''82.84'

$ws.Range('D17').Value = '2.104.60'
$ws.Range('E17').Value = '  -4.85%  '

$ws.Range('D18').Value = '''6.118'
$ws.Range('E18').Value = '  -1.95%  '

$ws.Range('D19').Value = '29.413.94'
$ws.Range('E19').Value = '  -0.65%  '

$ws.Range('D20').Value = '''227.53'
$ws.Range('E20').Value = '  -2.68%  '

$ws.Range('E21').Value = '  -0.43%  '

$ws.Range('D22').Value = '''0.9997'
$ws.Range('E22').Value = '  +0.00%  '

$ws.Range('D23').Value = '''7.538'
$ws.Range('E23').Value = '  -1.56%  '

$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').Value = '''3.951'
$ws.Range('E24').Value = '  -0.71%  '

$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').Value = '''1.001'
$ws.Range('E25').Value = '  +0.05%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''157.01'
$ws.Range('E26').Value = '  +1.52%  '

$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '''0.1396'
$ws.Range('E27').Value = '  +0.33%  '

$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''8.355'
$ws.Range('E28').Value = '  -1.14%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''17.64'
$ws.Range('E29').Value = '  -0.38%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''1.465'
$ws.Range('E30').Value = '  -0.99%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.05680'
$ws.Range('E31').Value = '  -2.85%  '

$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '''1.252'
$ws.Range('E32').Value = '  -0.22%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''4.121'
$ws.Range('E33').Value = '  +0.38%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '''4.022'
$ws.Range('E34').Value = '  -0.38%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '''1.843'
$ws.Range('E35').Value = '  -2.08%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '''1.154'
$ws.Range('E36').Value = '  -1.28%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''0.7153'
$ws.Range('E37').Value = '  -0.80%  '

$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '''2.596'
$ws.Range('E38').Value = '  +0.35%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.260.76'
$ws.Range('E39').Value = '  +1.30%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.01814'
$ws.Range('E40').Value = '  +1.67%  '

$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '''2.779'
$ws.Range('E41').Value = '  -0.70%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''6.222'
$ws.Range('E42').Value = '  +0.81%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''0.9105'
$ws.Range('E43').Value = '  +0.42%  '

$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '''0.9996'
$ws.Range('E44').Value = '  +0.00%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '''101.16'
$ws.Range('E45').Value = '  -0.92%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''66.15'
$ws.Range('E46').Value = '  -1.63%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '''7.055'
$ws.Range('E47').Value = '  -3.73%  '

$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '''0.00000000117'
$ws.Range('E48').Value = '  -0.69%  '

$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').Value = '''0.4036'
$ws.Range('E49').Value = '  -0.20%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''9.074'
$ws.Range('E50').Value = '  -0.87%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '''1.682'
$ws.Range('E51').Value = '  -1.57%  '

$ws.Range('D4').Style = "Normal"
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
